$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39. This shifts the former rows 39..109
# down to 40..110, preserving all of their data and formatting.
$ws.Rows.Item(39).Insert()

# The newly inserted row 39 is blank. Use the row right below it (which
# now holds what used to be row 39) as a formatting/content template, by
# copying it into the new row 39. We will then overwrite just the cells
# that actually differ (D, J, K, L, M, P) to create the new data record.
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(39).PasteSpecial()

# Now set the new record's distinct values on row 39.
$ws.Cells.Item(39, 4).Value  = 44791   # D39 Fecha
$ws.Cells.Item(39, 10).Value = 130     # J39 Volumen
$ws.Cells.Item(39, 11).Value = 6000    # K39 Precio minimo
$ws.Cells.Item(39, 12).Value = 10000   # L39 Precio maximo
$ws.Cells.Item(39, 13).Value = 8462    # M39 Precio promedio ponderado
$ws.Cells.Item(39, 16).Value = 141     # P39 Precio $/Kg
